$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 12.87519999999999
$ws.Range("E6").Value = 12.04
$ws.Range("D7").Value = -6.926799999999997
$ws.Range("A8").Value = -21.20020000000001
$ws.Range("E9").Value = 9.720899999999988
$ws.Range("A10").Value = -20.51829999999997
$ws.Range("E10").Value = 11.6079
$ws.Range("A12").Value = -22.65440000000004
$ws.Range("C13").Value = -12.55689999999999
$ws.Range("A18").Value = -22.45740000000004
$ws.Range("D20").Value = -8.541899999999995
